$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing text interpretation so that
# numeric-looking strings (e.g. "1.00", "0.720") are preserved exactly,
# then restore the cell to its original (unstyled) Normal style so no
# stray formatting is introduced.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "68.897.27"
Set-TextValue $ws.Range("E2") "  -3.70%  "

# Row 3
Set-TextValue $ws.Range("D3") "3.709.90"
Set-TextValue $ws.Range("E3") "  -4.34%  "

# Row 4
Set-TextValue $ws.Range("E4") "  -0.01%  "

# Row 5
Set-TextValue $ws.Range("D5") "599.54"
Set-TextValue $ws.Range("E5") "  -0.73%  "

# Row 6
Set-TextValue $ws.Range("D6") "181.60"
Set-TextValue $ws.Range("E6") "  +4.87%  "

# Row 7
Set-TextValue $ws.Range("D7") "3.707.02"
Set-TextValue $ws.Range("E7") "  -4.32%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.632"
Set-TextValue $ws.Range("E8") "  -5.67%  "

# Row 9
Set-TextValue $ws.Range("D9") "1.00"
Set-TextValue $ws.Range("E9") "  +0.02%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.720"
Set-TextValue $ws.Range("E10") "  -4.01%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.162"
Set-TextValue $ws.Range("E11") "  -8.98%  "

# Row 12
Set-TextValue $ws.Range("D12") "56.35"
Set-TextValue $ws.Range("E12") "  +4.24%  "

# Row 13
Set-TextValue $ws.Range("D13") "0.0000296"
Set-TextValue $ws.Range("E13") "  -8.52%  "

# Row 14
Set-TextValue $ws.Range("D14") "10.71"
Set-TextValue $ws.Range("E14") "  -7.33%  "

# Row 15
Set-TextValue $ws.Range("D15") "4.287.77"
Set-TextValue $ws.Range("E15") "  -4.69%  "

# Row 16
Set-TextValue $ws.Range("D16") "3.697.37"
Set-TextValue $ws.Range("E16") "  -4.66%  "

# Row 17
Set-TextValue $ws.Range("D17") "19.49"
Set-TextValue $ws.Range("E17") "  -8.00%  "

# Row 18
Set-TextValue $ws.Range("E18") "  -1.94%  "

# Row 19
Set-TextValue $ws.Range("D19") "12.96"
Set-TextValue $ws.Range("E19") "  -7.39%  "

# Row 20
Set-TextValue $ws.Range("E20") "  -7.01%  "

# Row 21
Set-TextValue $ws.Range("D21") "68.601.58"
Set-TextValue $ws.Range("E21") "  -3.76%  "

# Row 22
Set-TextValue $ws.Range("D22") "414.10"
Set-TextValue $ws.Range("E22") "  -6.12%  "

# Row 23
Set-TextValue $ws.Range("D23") "4.65"
Set-TextValue $ws.Range("E23") "  -2.84%  "

# Row 24
Set-TextValue $ws.Range("D24") "89.16"
Set-TextValue $ws.Range("E24") "  -5.77%  "

# Row 25
Set-TextValue $ws.Range("E25") "  -7.91%  "

# Row 26
Set-TextValue $ws.Range("D26") "12.76"
Set-TextValue $ws.Range("E26") "  -8.32%  "

# Row 27
Set-TextValue $ws.Range("D27") "10.95"
Set-TextValue $ws.Range("E27") "  -6.28%  "

# Row 28
Set-TextValue $ws.Range("E28") "  -3.62%  "

# Row 29
Set-TextValue $ws.Range("D29") "6.09"
Set-TextValue $ws.Range("E29") "  +2.07%  "

# Row 30
Set-TextValue $ws.Range("D30") "9.67"
Set-TextValue $ws.Range("E30") "  -7.72%  "

# Row 31
Set-TextValue $ws.Range("D31") "33.10"
Set-TextValue $ws.Range("E31") "  -6.19%  "

# Row 32
Set-TextValue $ws.Range("D32") "7.38"
Set-TextValue $ws.Range("E32") "  -13.93%  "

# Row 33
Set-TextValue $ws.Range("D33") "12.51"
Set-TextValue $ws.Range("E33") "  -7.90%  "

# Row 34
Set-TextValue $ws.Range("E34") "  -5.11%  "

# Row 35
Set-TextValue $ws.Range("D35") "43.90"
Set-TextValue $ws.Range("E35") "  -8.42%  "

# Row 36
Set-TextValue $ws.Range("D36") "64.97"
Set-TextValue $ws.Range("E36") "  -6.92%  "

# Row 37
Set-TextValue $ws.Range("D37") "605.15"
Set-TextValue $ws.Range("E37") "  -4.19%  "

# Row 38
Set-TextValue $ws.Range("D38") "0.0₃0886"
Set-TextValue $ws.Range("E38") "  -11.61%  "

# Row 39
Set-TextValue $ws.Range("D39") "0.406"
Set-TextValue $ws.Range("E39") "  -7.29%  "

# Row 40
Set-TextValue $ws.Range("E40") "  +0.12%  "

# Row 41
Set-TextValue $ws.Range("D41") "0.999"
Set-TextValue $ws.Range("E41") "  -0.11%  "

# Row 42
Set-TextValue $ws.Range("D42") "0.139"
Set-TextValue $ws.Range("E42") "  -5.86%  "

# Row 43
Set-TextValue $ws.Range("E43") "  -6.67%  "

# Row 44
Set-TextValue $ws.Range("D44") "0.0444"
Set-TextValue $ws.Range("E44") "  -6.20%  "

# Row 45
Set-TextValue $ws.Range("D45") "2.68"
Set-TextValue $ws.Range("E45") "  -6.51%  "

# Row 46
Set-TextValue $ws.Range("D46") "2.79"
Set-TextValue $ws.Range("E46") "  -11.86%  "

# Row 47
Set-TextValue $ws.Range("D47") "9.15"
Set-TextValue $ws.Range("E47") "  -10.61%  "

# Row 48
Set-TextValue $ws.Range("B48") "WEMIXToken"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D48") "2.74"
Set-TextValue $ws.Range("E48") "  -6.74%  "

# Row 49
Set-TextValue $ws.Range("B49") "Stellar"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D49") "0.136"
Set-TextValue $ws.Range("E49") "  -6.25%  "

# Row 50
Set-TextValue $ws.Range("D50") "2.755.38"
Set-TextValue $ws.Range("E50") "  -5.01%  "

# Row 51
Set-TextValue $ws.Range("D51") "3.08"
Set-TextValue $ws.Range("E51") "  -4.93%  "
